$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (13) and last column (M) to shrink the 13x13 matrix to 12x12
$ws.Rows(13).Delete()
$ws.Columns("M").Delete()

# Shift the judge labels in the header row and first column so that the
# "lena" label disappears and "yzxn" is appended at the end (matches the
# target workbook exactly).
$labels = @("cvai","eoce","jkwa","jvfs","lant","lskw","pdav","qoth","szch","xgju","yzxn")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $col = $i + 2          # column B is index 2
    $row = $i + 2           # row 2 is first data row
    $ws.Cells.Item(1, $col).Value = $labels[$i]
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
